# Apply "Ajout draft mapping f595a2bd5e53be80aa00972cfd76eee4a5f7087b"
#
# 1. Update the IG "Date" metadata value on the Metadata sheet.
# 2. Add a new mapping column ("Mapping: Spécification métier vers l'extension
#    ROR OrganizationCreationDate") to the Elements sheet, filled in only for
#    the Extension.value[x] row with "dateCreation".

$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the Date property value -----------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- 2. Elements sheet: append the new mapping column -----------------------
$elements = $wb.Worksheets.Item("Elements")

# Column AL is column 38 (A=1 ... Z=26, AA=27 ... AK=37, AL=38)
$newCol = 38
$lastExistingCol = 37

# Header (row 1) - copy formatting from the existing last header cell (AK1)
# and then set the new header text.
$elements.Cells.Item(1, $lastExistingCol).Copy($elements.Cells.Item(1, $newCol))
$elements.Cells.Item(1, $newCol).Value = "Mapping: Spécification métier vers l'extension ROR OrganizationCreationDate"

# Data rows (2-6) - copy formatting from the neighboring "Mapping: RIM Mapping"
# column cell so the new cells share the same body style, then set values.
# Rows 2-5 stay blank, row 6 (Extension.value[x]) gets "dateCreation".
for ($row = 2; $row -le 6; $row++) {
    $elements.Cells.Item($row, $lastExistingCol).Copy($elements.Cells.Item($row, $newCol))
    $elements.Cells.Item($row, $newCol).Value = ""
}
$elements.Cells.Item(6, $newCol).Value = "dateCreation"

# Column width for the new column (matches the authored width of ~82.04)
$elements.Columns.Item($newCol).ColumnWidth = 81.15
